$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) with full party names
$ws.Range("B1").Value = "AS - Åland Coalition (Åländsk Samling, AS)"
$ws.Range("C1").Value = "KD - Christian Democrats (Kristillisdemokraatit, KD), known until 25 May 2001 as Christian League (SKL, Suomen Kristillinen Liitto)"
$ws.Range("D1").Value = "KESK - Centre Party (Suomen Keskusta, KESK)"
$ws.Range("E1").Value = "KOK - National Coalition (Kansallinen Kokoomus, KOK)"
$ws.Range("F1").Value = "LKP - Liberal People’s Party (Liberaalinen Kansanpuolue, LKP)"
$ws.Range("G1").Value = "PS - The Finns Party (Perussuomalaiset, PS), known until  as True Finns (PS, Perussuomalaiset/ Sannfinländarna), known from  until 13 October 1995 as Finnish Rural Party (SMP, Suomen maaseudun puolue)"
$ws.Range("H1").Value = "SDP - Social Democratic Party (Sosialidemokraattinen Puolue, SDP)"
$ws.Range("I1").Value = "SFP - Swedish People’s Party (Svenska Folkpartiet, SFP)"
$ws.Range("J1").Value = "VAS - Left-Wing Alliance (Vasemmistoliitto , VAS)"
$ws.Range("K1").Value = "VIHR - Green League (Vihrea Liitto, VIHR)"
$ws.Range("L1").Value = "EKO - Ecological Party (Ekologinen Puolue, EKO)"
$ws.Range("M1").Value = "NUORS - Progressive Finnish Party (Nuorsuomalainen Puolue, NUORS)"
$ws.Range("N1").Value = "Other - Other (-, Other)"

# Fix floating point rounding artifacts in rows 8-9
$ws.Range("C8").Value = 5
$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 9
$ws.Range("K8").Value = 12

$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 17
$ws.Range("E9").Value = 22
$ws.Range("G9").Value = 23
$ws.Range("H9").Value = 24
$ws.Range("I9").Value = 7
$ws.Range("J9").Value = 10
$ws.Range("K9").Value = 13
